# Applies the 2026-01-09 odds refresh described in the commit diff:
#  - updates ~125 existing odds cells across rows 2-24
#  - appends a new fixture row 25 (Mazatlan FC vs FC Juarez)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing odds cells per diff ---
# Row 2
$ws.Cells.Item(2, 7).Value = 3.3
$ws.Cells.Item(2, 12).Value = 1.42
$ws.Cells.Item(2, 13).Value = 1.08
$ws.Cells.Item(2, 14).Value = 3.7
$ws.Cells.Item(2, 15).Value = 1.34
$ws.Cells.Item(2, 16).Value = 1.91
$ws.Cells.Item(2, 17).Value = 2.04
$ws.Cells.Item(2, 18).Value = 1.35
$ws.Cells.Item(2, 19).Value = 3.6
$ws.Cells.Item(2, 20).Value = 1.8
$ws.Cells.Item(2, 21).Value = 2.16
$ws.Cells.Item(2, 23).Value = 1.43
$ws.Cells.Item(2, 24).Value = 12
$ws.Cells.Item(2, 31).Value = 26
$ws.Cells.Item(2, 40).Value = 36

# Row 3
$ws.Cells.Item(3, 13).Value = 1.05
$ws.Cells.Item(3, 15).Value = 1.24
$ws.Cells.Item(3, 20).Value = 1.71

# Row 4
$ws.Cells.Item(4, 6).Value = 1.84
$ws.Cells.Item(4, 9).Value = 4.8
$ws.Cells.Item(4, 20).Value = 1.71
$ws.Cells.Item(4, 26).Value = 36
$ws.Cells.Item(4, 27).Value = 110
$ws.Cells.Item(4, 31).Value = 55
$ws.Cells.Item(4, 35).Value = 60
$ws.Cells.Item(4, 37).Value = 21
$ws.Cells.Item(4, 41).Value = 55

# Row 5
$ws.Cells.Item(5, 6).Value = 1.7
$ws.Cells.Item(5, 9).Value = 7.6
$ws.Cells.Item(5, 10).Value = 3.4
$ws.Cells.Item(5, 11).Value = 3.65
$ws.Cells.Item(5, 14).Value = 2.56
$ws.Cells.Item(5, 17).Value = 2.56
$ws.Cells.Item(5, 20).Value = 2.32
$ws.Cells.Item(5, 23).Value = 2.24
$ws.Cells.Item(5, 24).Value = 9.2
$ws.Cells.Item(5, 28).Value = 6
$ws.Cells.Item(5, 29).Value = 8.8
$ws.Cells.Item(5, 30).Value = 34
$ws.Cells.Item(5, 32).Value = 9
$ws.Cells.Item(5, 34).Value = 38
$ws.Cells.Item(5, 35).Value = 200
$ws.Cells.Item(5, 39).Value = 340
$ws.Cells.Item(5, 41).Value = 370

# Row 6
$ws.Cells.Item(6, 6).Value = 2.48
$ws.Cells.Item(6, 7).Value = 2.76
$ws.Cells.Item(6, 9).Value = 3.9
$ws.Cells.Item(6, 14).Value = 2.28
$ws.Cells.Item(6, 16).Value = 1.43
$ws.Cells.Item(6, 17).Value = 2.8
$ws.Cells.Item(6, 22).Value = 1.35
$ws.Cells.Item(6, 23).Value = 1.56
$ws.Cells.Item(6, 24).Value = 7.8
$ws.Cells.Item(6, 33).Value = 14.5

# Row 7
$ws.Cells.Item(7, 9).Value = 2.3
$ws.Cells.Item(7, 11).Value = 3.4
$ws.Cells.Item(7, 12).Value = 1.65
$ws.Cells.Item(7, 17).Value = 3.1
$ws.Cells.Item(7, 18).Value = 1.12
$ws.Cells.Item(7, 20).Value = 2.44
$ws.Cells.Item(7, 29).Value = 1000

# Row 8
$ws.Cells.Item(8, 7).Value = 2.06
$ws.Cells.Item(8, 10).Value = 3.4
$ws.Cells.Item(8, 12).Value = 1.37
$ws.Cells.Item(8, 15).Value = 1.35
$ws.Cells.Item(8, 16).Value = 1.78
$ws.Cells.Item(8, 17).Value = 2.02
$ws.Cells.Item(8, 19).Value = 3.7
$ws.Cells.Item(8, 20).Value = 1.85
$ws.Cells.Item(8, 23).Value = 1.94

# Row 9
$ws.Cells.Item(9, 6).Value = 1.95
$ws.Cells.Item(9, 10).Value = 2.8
$ws.Cells.Item(9, 11).Value = 3.6
$ws.Cells.Item(9, 12).Value = 1.56
$ws.Cells.Item(9, 14).Value = 2.18
$ws.Cells.Item(9, 20).Value = 2.48
$ws.Cells.Item(9, 22).Value = 1.21

# Row 10
$ws.Cells.Item(10, 8).Value = 1.64
$ws.Cells.Item(10, 11).Value = 4.5
$ws.Cells.Item(10, 17).Value = 1.69
$ws.Cells.Item(10, 19).Value = 3

# Row 11
$ws.Cells.Item(11, 17).Value = 1.83

# Row 12
$ws.Cells.Item(12, 7).Value = 10
$ws.Cells.Item(12, 9).Value = 1.51
$ws.Cells.Item(12, 12).Value = 1.25
$ws.Cells.Item(12, 14).Value = 5
$ws.Cells.Item(12, 16).Value = 2.4
$ws.Cells.Item(12, 17).Value = 1.59
$ws.Cells.Item(12, 18).Value = 1.56
$ws.Cells.Item(12, 20).Value = 1.8
$ws.Cells.Item(12, 21).Value = 2.02
$ws.Cells.Item(12, 22).Value = 2.96
$ws.Cells.Item(12, 23).Value = 1.13
$ws.Cells.Item(12, 33).Value = 980
$ws.Cells.Item(12, 35).Value = 980

# Row 20
$ws.Cells.Item(20, 6).Value = 1.98
$ws.Cells.Item(20, 7).Value = 2.2
$ws.Cells.Item(20, 10).Value = 2.98
$ws.Cells.Item(20, 14).Value = 3.05
$ws.Cells.Item(20, 15).Value = 1.4
$ws.Cells.Item(20, 16).Value = 1.77
$ws.Cells.Item(20, 17).Value = 2.06
$ws.Cells.Item(20, 18).Value = 1.29
$ws.Cells.Item(20, 20).Value = 1.87
$ws.Cells.Item(20, 21).Value = 1.92
$ws.Cells.Item(20, 23).Value = 1.84

# Row 21
$ws.Cells.Item(21, 6).Value = 2.28
$ws.Cells.Item(21, 7).Value = 2.66
$ws.Cells.Item(21, 9).Value = 4
$ws.Cells.Item(21, 10).Value = 2.94
$ws.Cells.Item(21, 11).Value = 3.35
$ws.Cells.Item(21, 15).Value = 1.01
$ws.Cells.Item(21, 16).Value = 1.25
$ws.Cells.Item(21, 17).Value = 2.26
$ws.Cells.Item(21, 19).Value = 1.01
$ws.Cells.Item(21, 22).Value = 1.33
$ws.Cells.Item(21, 23).Value = 1.63

# Row 22
$ws.Cells.Item(22, 27).Value = 65

# Row 23
$ws.Cells.Item(23, 9).Value = 2.42
$ws.Cells.Item(23, 14).Value = 3.55
$ws.Cells.Item(23, 16).Value = 1.89
$ws.Cells.Item(23, 18).Value = 1.34
$ws.Cells.Item(23, 21).Value = 1.95
$ws.Cells.Item(23, 22).Value = 1.7

# Row 24
$ws.Cells.Item(24, 17).Value = 3.3

# --- Append new row 25: Mexican Liga MX, Mazatlan FC vs FC Juarez ---
$newRow = 25

# Text columns (League/Home/Away) - safe to assign directly
$ws.Cells.Item($newRow, 1).Value = "Mexican Liga MX"
$ws.Cells.Item($newRow, 4).Value = "Mazatlan FC"
$ws.Cells.Item($newRow, 5).Value = "FC Juarez"

# Date/Time columns must be forced to text format first so the
# COM layer stores the literal string instead of auto-converting
# "2026-01-09" / "22:00:00" into a date/time serial number.
$dateCell = $ws.Cells.Item($newRow, 2)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026-01-09"

$timeCell = $ws.Cells.Item($newRow, 3)
$timeCell.NumberFormat = "@"
$timeCell.Value = "22:00:00"

# Numeric odds columns (F .. AO)
$ws.Cells.Item($newRow, 6).Value = 3.05
$ws.Cells.Item($newRow, 7).Value = 3.45
$ws.Cells.Item($newRow, 8).Value = 2.38
$ws.Cells.Item($newRow, 9).Value = 2.6
$ws.Cells.Item($newRow, 10).Value = 3.3
$ws.Cells.Item($newRow, 11).Value = 3.7
$ws.Cells.Item($newRow, 12).Value = 1.01
$ws.Cells.Item($newRow, 13).Value = 1.01
$ws.Cells.Item($newRow, 14).Value = 1.84
$ws.Cells.Item($newRow, 15).Value = 1.34
$ws.Cells.Item($newRow, 16).Value = 1.84
$ws.Cells.Item($newRow, 17).Value = 2
$ws.Cells.Item($newRow, 18).Value = 1.28
$ws.Cells.Item($newRow, 19).Value = 3.15
$ws.Cells.Item($newRow, 20).Value = 1.64
$ws.Cells.Item($newRow, 21).Value = 1.78
$ws.Cells.Item($newRow, 22).Value = 1.62
$ws.Cells.Item($newRow, 23).Value = 1.4
$ws.Cells.Item($newRow, 24).Value = 19
$ws.Cells.Item($newRow, 25).Value = 12.5
$ws.Cells.Item($newRow, 26).Value = 23
$ws.Cells.Item($newRow, 27).Value = 44
$ws.Cells.Item($newRow, 28).Value = 15
$ws.Cells.Item($newRow, 29).Value = 9.4
$ws.Cells.Item($newRow, 30).Value = 17.5
$ws.Cells.Item($newRow, 31).Value = 34
$ws.Cells.Item($newRow, 32).Value = 30
$ws.Cells.Item($newRow, 33).Value = 20
$ws.Cells.Item($newRow, 34).Value = 26
$ws.Cells.Item($newRow, 35).Value = 65
$ws.Cells.Item($newRow, 36).Value = 70
$ws.Cells.Item($newRow, 37).Value = 50
$ws.Cells.Item($newRow, 38).Value = 65
$ws.Cells.Item($newRow, 39).Value = 1000
$ws.Cells.Item($newRow, 40).Value = 1000
$ws.Cells.Item($newRow, 41).Value = 1000
